# Add two new weekly price rows (15 and 16) for "Femacal de La Calera" /
# Chirimoya, matching the layout of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 -----------------------------------------------------------
$ws.Cells.Item(15, 1).Value  = 3
$ws.Cells.Item(15, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(15, 3).Value  = "Coquimbo"
$ws.Cells.Item(15, 4).Value  = 44448
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value  = 5
$ws.Cells.Item(15, 6).Value  = "Fruta"
$ws.Cells.Item(15, 7).Value  = 100107
$ws.Cells.Item(15, 8).Value  = "Otros"
$ws.Cells.Item(15, 9).Value  = 100107002
$ws.Cells.Item(15, 10).Value = "Chirimoya"
$ws.Cells.Item(15, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 45
$ws.Cells.Item(15, 14).Value = 30000
$ws.Cells.Item(15, 15).Value = 30000
$ws.Cells.Item(15, 16).Value = 30000
$ws.Cells.Item(15, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(15, 19).Value = 3000
$ws.Cells.Item(15, 20).Value = 10

# --- Row 16 -----------------------------------------------------------
$ws.Cells.Item(16, 1).Value  = 3
$ws.Cells.Item(16, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(16, 3).Value  = "Coquimbo"
$ws.Cells.Item(16, 4).Value  = 44448
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value  = 5
$ws.Cells.Item(16, 6).Value  = "Fruta"
$ws.Cells.Item(16, 7).Value  = 100107
$ws.Cells.Item(16, 8).Value  = "Otros"
$ws.Cells.Item(16, 9).Value  = 100107002
$ws.Cells.Item(16, 10).Value = "Chirimoya"
$ws.Cells.Item(16, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 40
$ws.Cells.Item(16, 14).Value = 27000
$ws.Cells.Item(16, 15).Value = 27000
$ws.Cells.Item(16, 16).Value = 27000
$ws.Cells.Item(16, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(16, 19).Value = 2700
$ws.Cells.Item(16, 20).Value = 10
